$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new text value (all target cells are Text-typed, like the originals)
$updates = @{
    'D2' = '277.61'
    'E2' = '2.04%'
    'G2' = '22'
    'D3' = '27.29'
    'E3' = '1.67%'
    'G3' = '22'
    'D4' = '4.943'
    'E4' = '0.76%'
    'G4' = '22'
    'D5' = '0.06404'
    'E5' = '1.47%'
    'G5' = '22'
    'D6' = '6.979'
    'E6' = '1.20%'
    'G6' = '22'
    'D7' = '1.260'
    'E7' = '-9.05%'
    'G7' = '22'
    'D8' = '0.8810'
    'E8' = '-0.28%'
    'G8' = '22'
    'D9' = '0.1523'
    'E9' = '3.85%'
    'G9' = '22'
    'D10' = '0.05119'
    'E10' = '0.52%'
    'G10' = '22'
    'D11' = '0.07511'
    'E11' = '1.51%'
    'G11' = '22'
    'D12' = '0.02960'
    'E12' = '-6.63%'
    'G12' = '22'
    'D13' = '0.09016'
    'E13' = '-0.26%'
    'G13' = '22'
    'D14' = '0.001564'
    'E14' = '0.04%'
    'G14' = '22'
    'D15' = '0.0006402'
    'E15' = '1.25%'
    'G15' = '22'
    'D16' = '0.005910'
    'E16' = '-2.53%'
    'G16' = '22'
    'D17' = '3.457'
    'E17' = '-0.34%'
    'G17' = '22'
    'D18' = '3.322'
    'E18' = '-0.96%'
    'G18' = '22'
    'D19' = '2.285'
    'E19' = '0.04%'
    'G19' = '22'
    'E20' = '1.36%'
    'G20' = '22'
    'D21' = '0.1338'
    'E21' = '0.22%'
    'G21' = '22'
    'D22' = '3.908'
    'E22' = '-0.07%'
    'G22' = '22'
    'D23' = '0.04433'
    'E23' = '2.17%'
    'G23' = '22'
    'D24' = '0.001172'
    'E24' = '-0.47%'
    'G24' = '22'
    'D25' = '0.003873'
    'E25' = '6.36%'
    'G25' = '22'
    'D26' = '0.0001200'
    'E26' = '-0.11%'
    'G26' = '22'
    'E27' = '13.91%'
    'G27' = '22'
    'G28' = '22'
    'G29' = '22'
    'G30' = '22'
    'G31' = '22'
    'G32' = '22'
    'G33' = '22'
    'G34' = '22'
    'G35' = '22'
    'G36' = '22'
    'G37' = '22'
    'G38' = '22'
    'G39' = '22'
    'E40' = '2.79%'
    'G40' = '22'
    'D41' = '0.006804'
    'E41' = '2.80%'
    'G41' = '22'
    'E42' = '1.41%'
    'G42' = '22'
    'D43' = '0.002020'
    'E43' = '-6.58%'
    'G43' = '22'
    'D44' = '0.01122'
    'E44' = '-10.95%'
    'G44' = '22'
    'D45' = '0.00005178'
    'E45' = '-2.97%'
    'G45' = '22'
    'G46' = '22'
    'D47' = '0.02026'
    'E47' = '-4.57%'
    'G47' = '22'
    'G48' = '22'
    'G49' = '22'
    'G50' = '22'
    'G51' = '22'
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
    $cell.Style = "Normal"
}
